# TN ALL STAR prep: unhide helper columns, add name/location formulas,
# update the active selection, and lift sheet protection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The helper columns (F:H) are protected/hidden; unprotect first so we can
# edit them and change their visibility.
$ws.Unprotect()

# Add the concatenation helper formulas that feed the scorebug graphics.
$ws.Range("F7").Formula = '=B7 & ", " & B8'
$ws.Range("F9").Formula = '=B2 & " " & C2'
$ws.Range("F10").Formula = '=B3 & " " & C3'

# Unhide the helper columns F:H. Column G previously had no real stored
# width (it was hidden at width 0), so give it a sensible visible width;
# F and H already carry their real widths, so just toggle visibility.
$ws.Columns.Item(7).ColumnWidth = 8.33
$ws.Columns.Item(6).Hidden = $false
$ws.Columns.Item(7).Hidden = $false
$ws.Columns.Item(8).Hidden = $false

# Move/extend the active selection to E9:E10.
$ws.Range("E9:E10").Select() | Out-Null
